$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A31").Value = "resumeUrl"
Write-Host "done"
